# Updated cryptos list data (Price + Volume(1h) columns, plus two swapped rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '41.670.30'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -1.22%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.173.92'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -2.75%  '

# Row 4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '238.47'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.96%  '

# Row 6
$ws.Range("E6").Value = '  -2.30%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '72.49'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -2.61%  '

# Row 8
$ws.Range("E8").Value = '  -0.13%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.582'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -3.15%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '40.55'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -4.12%  '

# Row 11
$ws.Range("E11").Value = '  -4.85%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '54.62'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -3.36%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '6.74'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -2.91%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.100'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -3.10%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '2.501.29'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.72%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '14.40'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +0.32%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '2.170.74'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -2.56%  '

# Row 18
$ws.Range("E18").Value = '  -6.71%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '41.651.26'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.08%  '

# Row 20
$ws.Range("E20").Value = '  -1.99%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '70.18'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -3.64%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '5.80'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -6.80%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '10.05'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -13.07%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '226.38'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.63%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.04'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.40%  '

# Row 26
$ws.Range("E26").Value = '  +0.16%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '10.77'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -5.21%  '

# Row 28
$ws.Range("E28").Value = '  -10.02%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '2.20'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -3.23%  '

# Row 30
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.18'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -1.14%  '

# Row 31
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '170.79'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.98%  '

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '19.86'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -3.69%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '32.70'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +9.54%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.0777'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -3.31%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '5.31'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -5.91%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.120'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -3.45%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '4.34'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.26%  '

# Row 38
$ws.Range("E38").Value = '  -6.05%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.0312'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.75%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '12.12'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -7.86%  '

# Row 41
$ws.Range("E41").Value = '  -1.61%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.37'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -5.93%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '59.26'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -8.56%  '

# Row 44
$ws.Range("E44").Value = '  -4.60%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '8.46'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -3.04%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0964'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -3.60%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '97.46'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -6.80%  '

# Row 48
$ws.Range("E48").Value = '  -4.33%  '

# Row 49
$ws.Range("E49").Value = '  -4.86%  '

# Row 50
$ws.Range("E50").Value = '  -6.24%  '

# Row 51
$ws.Range("E51").Value = '  -2.09%  '
